$p = $ppt.ActivePresentation

# Slide 20 ("Sample File") -> shape 2 ("Content Placeholder 2") holds the
# PF rule listing in the Courier-formatted paragraphs.
$s = $p.Slides.Item(20)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 4 is "pass in on em0 inet proto tcp all"
$para = $tr.Paragraphs(4, 1)

if ($para.Text.IndexOf("pass in on em0 ") -ne 0) {
    throw "Unexpected paragraph text: $($para.Text)"
}

# Within that paragraph, characters 9-11 are "on " (inside the leading
# "pass in on em0 " run). Replacing just that slice splits the run into
# three runs - "pass in ", "log on ", "em0 " - while leaving the
# trailing "inet" / " proto " / "tcp" / " all" runs untouched.
$target = $para.Characters(9, 3)
if ($target.Text -ne "on ") {
    throw "Unexpected characters slice: [$($target.Text)]"
}
$target.Text = "log on "
